$d = $word.ActiveDocument

# --- 1. Insert a new ListParagraph bullet after "Den mellemste sværhedsgrad ... sættes op."
#        and before the "Ændringer foretaget i programmet" (Heading2) paragraph. ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("bør vi derfor overveje om skal sættes op.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Collapse(0)
    $r1.InsertAfter([char]13 + "Hvis børnene fortrød at ville placere et tal, brugte de ”Fjern tal”-knappen til at lukke dialogen med.")
}

# --- 2. Insert a new plain paragraph after "Desuden ændrede vi maskottens tekst ... feltet kun har én løsning."
#        and before the "Funktionstester" (Heading1) paragraph. ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("angives at feltet kun har én løsning.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)
    $r2.InsertAfter([char]13 + "Vi tilføjede også en ”Luk”-knap til nummervælgerdialogen, så ”Fjern tal” forhåbentlig vil blive brugt korrekt.")
}
